# Backup QR Scanner data - append newly scanned log rows to the Histology
# scan log sheet (rows 28-63), matching columns:
#   A: Student ID, B: Subject, C: Log Date, D: Log Time, E: Type, F: User
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for new rows 28-63 (Student ID, Subject, Log Date, Log Time, Type, User)
$data = @(
    ,@("241220", "Histology", "17/12/2025", "08:42:30", "Scan", "norhan.mohamed@med.asu.edu.eg")
    ,@("241156", "Histology", "17/12/2025", "08:42:33", "Scan", "norhan.mohamed@med.asu.edu.eg")
    ,@("241160", "Histology", "17/12/2025", "08:42:37", "Scan", "norhan.mohamed@med.asu.edu.eg")
    ,@("241212", "Histology", "17/12/2025", "08:42:40", "Scan", "norhan.mohamed@med.asu.edu.eg")
    ,@("241196", "Histology", "17/12/2025", "08:42:43", "Scan", "norhan.mohamed@med.asu.edu.eg")
    ,@("241221", "Histology", "17/12/2025", "08:42:46", "Scan", "norhan.mohamed@med.asu.edu.eg")
    ,@("241131", "Histology", "17/12/2025", "08:42:51", "Scan", "norhan.mohamed@med.asu.edu.eg")
    ,@("241158", "Histology", "17/12/2025", "08:42:53", "Scan", "norhan.mohamed@med.asu.edu.eg")
    ,@("241113", "Histology", "17/12/2025", "08:43:01", "Scan", "norhan.mohamed@med.asu.edu.eg")
    ,@("241202", "Histology", "17/12/2025", "08:43:06", "Scan", "norhan.mohamed@med.asu.edu.eg")
    ,@("241190", "Histology", "17/12/2025", "08:43:11", "Scan", "norhan.mohamed@med.asu.edu.eg")
    ,@("241188", "Histology", "17/12/2025", "08:43:19", "Scan", "norhan.mohamed@med.asu.edu.eg")
    ,@("241209", "Histology", "17/12/2025", "08:43:25", "Scan", "norhan.mohamed@med.asu.edu.eg")
    ,@("241193", "Histology", "17/12/2025", "08:43:29", "Scan", "norhan.mohamed@med.asu.edu.eg")
    ,@("241195", "Histology", "17/12/2025", "08:43:32", "Scan", "norhan.mohamed@med.asu.edu.eg")
    ,@("241237", "Histology", "17/12/2025", "08:43:36", "Scan", "norhan.mohamed@med.asu.edu.eg")
    ,@("241199", "Histology", "17/12/2025", "08:43:39", "Scan", "norhan.mohamed@med.asu.edu.eg")
    ,@("241189", "Histology", "17/12/2025", "08:43:43", "Scan", "norhan.mohamed@med.asu.edu.eg")
    ,@("241205", "Histology", "17/12/2025", "08:43:47", "Scan", "norhan.mohamed@med.asu.edu.eg")
    ,@("241176", "Histology", "17/12/2025", "08:43:56", "Scan", "norhan.mohamed@med.asu.edu.eg")
    ,@("241167", "Histology", "17/12/2025", "08:43:58", "Scan", "norhan.mohamed@med.asu.edu.eg")
    ,@("241155", "Histology", "17/12/2025", "08:44:02", "Scan", "norhan.mohamed@med.asu.edu.eg")
    ,@("241126", "Histology", "17/12/2025", "08:44:06", "Scan", "norhan.mohamed@med.asu.edu.eg")
    ,@("241147", "Histology", "17/12/2025", "08:44:10", "Scan", "norhan.mohamed@med.asu.edu.eg")
    ,@("241187", "Histology", "17/12/2025", "08:44:15", "Scan", "norhan.mohamed@med.asu.edu.eg")
    ,@("241171", "Histology", "17/12/2025", "08:44:17", "Scan", "norhan.mohamed@med.asu.edu.eg")
    ,@("241121", "Histology", "17/12/2025", "08:44:22", "Scan", "norhan.mohamed@med.asu.edu.eg")
    ,@("241225", "Histology", "17/12/2025", "08:44:39", "Scan", "norhan.mohamed@med.asu.edu.eg")
    ,@("241186", "Histology", "17/12/2025", "08:44:55", "Scan", "norhan.mohamed@med.asu.edu.eg")
    ,@("241208", "Histology", "17/12/2025", "08:45:01", "Scan", "norhan.mohamed@med.asu.edu.eg")
    ,@("241174", "Histology", "17/12/2025", "08:45:05", "Scan", "norhan.mohamed@med.asu.edu.eg")
    ,@("241218", "Histology", "17/12/2025", "08:46:42", "Scan", "norhan.mohamed@med.asu.edu.eg")
    ,@("241102", "Histology", "17/12/2025", "08:46:46", "Scan", "norhan.mohamed@med.asu.edu.eg")
    ,@("241168", "Histology", "17/12/2025", "08:46:49", "Scan", "norhan.mohamed@med.asu.edu.eg")
    ,@("241166", "Histology", "17/12/2025", "08:46:53", "Scan", "norhan.mohamed@med.asu.edu.eg")
    ,@("241161", "Histology", "17/12/2025", "08:47:03", "Scan", "norhan.mohamed@med.asu.edu.eg")
)

$startRow = 28
$endRow = $startRow + $data.Count - 1

# Column A values look like plain integers (e.g. "241220"); format the
# destination range as Text first so Excel keeps them as strings instead
# of auto-converting to numbers (matches the existing rows 2-27 which are
# all stored as text).
$ws.Range("A" + $startRow + ":A" + $endRow).NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    $ws.Cells.Item($r, 1).Value = $rowVals[0]
    $ws.Cells.Item($r, 2).Value = $rowVals[1]
    $ws.Cells.Item($r, 3).Value = $rowVals[2]
    $ws.Cells.Item($r, 4).Value = $rowVals[3]
    $ws.Cells.Item($r, 5).Value = $rowVals[4]
    $ws.Cells.Item($r, 6).Value = $rowVals[5]
}